# Scheduled market-data refresh: updates currentAveragePrice / Leve profit
# columns (H:N) across the per-job Profits sheets.
$wb = $excel.ActiveWorkbook

function Set-Cells($SheetName, $Row, $Values, $Clear) {
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($col in $Values.Keys) {
        $ws.Range($col + $Row).Value = $Values[$col]
    }
    if ($Clear) {
        foreach ($col in $Clear) {
            $ws.Range($col + $Row).ClearContents()
        }
    }
}

# ALC
Set-Cells "ALC" 19 @{ H = 706.9; I = 693.7778; J = 717.63635; K = 693.7778; L = 717.63635; M = -518.7778; N = -1067.63635 } $null
Set-Cells "ALC" 40 @{ H = 5578254; I = 7813812.5; J = 2001360.4; K = 7813812.5; L = 2001360.4; M = -7813637.5; N = -2001710.4 } $null
Set-Cells "ALC" 43 @{ H = 27779204; I = 62501310; J = 1517.5; K = 62501310; L = 1517.5; M = -62501241; N = -1655.5 } $null
Set-Cells "ALC" 52 @{ H = 166669180; I = 1031.3334; J = 333337340; K = 3094.0002; L = 1000012020; M = -2934.0002; N = -1000012340 } $null
Set-Cells "ALC" 116 @{ H = 2344; J = 2828.2727; L = 2828.2727; N = -9712.2727 } $null
Set-Cells "ALC" 132 @{ H = 44126908; I = 45463330; K = 136389990; M = -136387460 } $null

# ARM
Set-Cells "ARM" 113 @{ H = 0; J = 0; L = 0 } @("N")

# BSM
Set-Cells "BSM" 134 @{ H = 4635647; I = 1829.5714; K = 5488.7142; M = -2953.7142 } $null

# CRP
Set-Cells "CRP" 58 @{ H = 90910300; I = 100000936; J = 4000; K = 100000936; L = 4000; M = -100000733; N = -4406 } $null
Set-Cells "CRP" 81 @{ H = 32240; J = 32240; L = 32240; N = -34236 } $null
Set-Cells "CRP" 82 @{ H = 21995; J = 21995; L = 21995; N = -22717 } $null
Set-Cells "CRP" 84 @{ H = 32240; J = 32240; L = 96720; N = -106704 } $null
Set-Cells "CRP" 85 @{ H = 21995; J = 21995; L = 21995; N = -24491 } $null
Set-Cells "CRP" 118 @{ H = 59800; J = 59800; L = 59800; N = -63114 } $null
Set-Cells "CRP" 132 @{ H = 33337692; I = 4510.2856; J = 111115110; K = 13530.8568; L = 333345330; M = -11000.8568; N = -333350390 } $null
Set-Cells "CRP" 134 @{ H = 1369.3334; I = 1170; J = 2166.6667; K = 3510; L = 6500.000100000001; M = -975; N = -11570.0001 } $null
Set-Cells "CRP" 136 @{ H = 90910300; I = 100000936; J = 4000; K = 300002808; L = 12000; M = -300000258; N = -17100 } $null

# CUL
Set-Cells "CUL" 59 @{ H = 1625; I = 500; K = 1500; M = -960 } $null

# GSM
Set-Cells "GSM" 105 @{ H = 59800; J = 59800; L = 59800; N = -66788 } $null
Set-Cells "GSM" 140 @{ H = 44400; J = 44400; L = 44400; N = -54760 } $null
Set-Cells "GSM" 141 @{ H = 29254.334; J = 29254.334; L = 29254.334; N = -39614.334 } $null

# LTW
Set-Cells "LTW" 22 @{ H = 2113537.8; I = 2535445.5; J = 4000; K = 2535445.5; L = 4000; M = -2535150.5; N = -4590 } $null
Set-Cells "LTW" 27 @{ H = 2113537.8; I = 2535445.5; J = 4000; K = 2535445.5; L = 4000; M = -2535338.5; N = -4214 } $null
Set-Cells "LTW" 55 @{ H = 58829556; I = 9166.454; J = 166666940; K = 9166.454; L = 166666940; M = -8993.454; N = -166667286 } $null
Set-Cells "LTW" 132 @{ H = 33621444; I = 63493370; J = 15531.125; K = 190480110; L = 46593.375; M = -190477580; N = -51653.375 } $null
Set-Cells "LTW" 136 @{ H = 238100240; I = 228576830; J = 250004500; K = 685730490; L = 750013500; M = -685727940; N = -750018600 } $null

# WVR
Set-Cells "WVR" 62 @{ H = 4176.6665; I = 3812; K = 3812; M = -3188 } $null
Set-Cells "WVR" 65 @{ H = 4176.6665; I = 3812; K = 19060; M = -15940 } $null
Set-Cells "WVR" 132 @{ H = 105118.73; I = 206261.6; J = 20833; K = 618784.8; L = 62499; M = -616254.8; N = -67559 } $null
Set-Cells "WVR" 136 @{ H = 7144590.5; I = 13890729; J = 1620; K = 41672187; L = 4860; M = -41669637; N = -9960 } $null
Set-Cells "WVR" 137 @{ H = 0; J = 0; L = 0 } @("N")
Set-Cells "WVR" 140 @{ H = 46583.75; J = 46583.75; L = 46583.75; N = -56943.75 } $null
